# Add "Source" info (Developed by RAPID Team) to rows 2-7 in column C,
# matching the existing value already present in rows 9-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceText = "Developed by RAPID Team"

foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = $sourceText
}

# C2 previously carried the wrap-text style (s="2") while empty; after
# populating it with the source text it reverts to the default/Normal style.
$ws.Range("C2").Style = "Normal"

# Update the active selection from A10 to C1
$ws.Range("C1").Select()
